$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.895.77'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.888.33'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7634'
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('E6').Value = '  -0.71%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3124'
$ws.Range('E8').Value = '  -0.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.64'
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07164'
$ws.Range('E10').Value = '  -3.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08505'
$ws.Range('E11').Value = '  +4.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7622'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.366'
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.864.13'
$ws.Range('E14').Value = '  +1.39%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.70'
$ws.Range('E15').Value = '  +1.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.139'
$ws.Range('E16').Value = '  -1.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.812.74'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.77'
$ws.Range('E18').Value = '  -1.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.62'
$ws.Range('E19').Value = '  -0.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007798'
$ws.Range('E20').Value = '  -0.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.994'
$ws.Range('E22').Value = '  -1.87%  '
$ws.Range('B23').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.102.99'
$ws.Range('E23').Value = '  -1.08%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1616'
$ws.Range('E25').Value = '  +2.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.404'
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '161.87'
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.75'
$ws.Range('E28').Value = '  -0.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.035'
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('E30').Value = '  +1.42%  '
$ws.Range('E31').Value = '  -0.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.489'
$ws.Range('E32').Value = '  -0.39%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05446'
$ws.Range('E34').Value = '  -2.81%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.243'
$ws.Range('E35').Value = '  -0.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7428'
$ws.Range('E36').Value = '  -2.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9991'
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.697'
$ws.Range('E38').Value = '  +1.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01946'
$ws.Range('E39').Value = '  +0.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.782'
$ws.Range('E40').Value = '  -0.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4465'
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.103.14'
$ws.Range('E42').Value = '  -4.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.064'
$ws.Range('E43').Value = '  +1.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.87'
$ws.Range('E44').Value = '  -2.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8525'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.82'
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.867'
$ws.Range('E48').Value = '  -2.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.610'
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.998'
$ws.Range('E50').Value = '  -5.26%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06077'
$ws.Range('E51').Value = '  +0.64%  '
